$wb = $excel.ActiveWorkbook

# Citywide Totals
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("C2").Value = 28
$ws.Range("D2").Value = 42
$ws.Range("D3").Value = 68
$ws.Range("J3").Value = 85
$ws.Range("F4").Value = 3
$ws.Range("B6").Value = 173
$ws.Range("C6").Value = 223
$ws.Range("D6").Value = 199
$ws.Range("E6").Value = 208
$ws.Range("F6").Value = 240
$ws.Range("G6").Value = 223
$ws.Range("H6").Value = 198
$ws.Range("I6").Value = 263
$ws.Range("J6").Value = 187
$ws.Range("B7").Value = 237
$ws.Range("C7").Value = 302
$ws.Range("D7").Value = 313
$ws.Range("E7").Value = 312
$ws.Range("F7").Value = 342
$ws.Range("G7").Value = 330
$ws.Range("H7").Value = 304
$ws.Range("I7").Value = 407
$ws.Range("J7").Value = 345

# By Neighborhood
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("F4").Value = 5
$ws.Range("C7").Value = 24
$ws.Range("F7").Value = 23
$ws.Range("H7").Value = 21
$ws.Range("I7").Value = 18
$ws.Range("J7").Value = 26
$ws.Range("C12").Value = 4
$ws.Range("D16").Value = 2
$ws.Range("D25").Value = 4
$ws.Range("B26").Value = 20
$ws.Range("G26").Value = 24
$ws.Range("D28").Value = 4
$ws.Range("E30").Value = 21
$ws.Range("F30").Value = 29
$ws.Range("C34").Value = 15
$ws.Range("E34").Value = 13
$ws.Range("I34").Value = 22
$ws.Range("J43").Value = 3
$ws.Range("D51").Value = 42
$ws.Range("E51").Value = 47
$ws.Range("H51").Value = 36
$ws.Range("D69").Value = 2
$ws.Range("D72").Value = 5
$ws.Range("C73").Value = 12
$ws.Range("C82").Value = 5
$ws.Range("G87").Value = 5
$ws.Range("B94").Value = 237
$ws.Range("C94").Value = 302
$ws.Range("D94").Value = 313
$ws.Range("E94").Value = 312
$ws.Range("F94").Value = 342
$ws.Range("G94").Value = 330
$ws.Range("H94").Value = 304
$ws.Range("I94").Value = 407
$ws.Range("J94").Value = 345

# Rogers Park
$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("D5").Value = 4
$ws.Range("D6").Value = 5

# Roseland
$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("C5").Value = 8
$ws.Range("C6").Value = 12

# Gage Park
$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("D4").Value = 4
$ws.Range("D5").Value = 4

# Austin
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("C2").Value = 3
$ws.Range("J3").Value = 5
$ws.Range("F5").Value = 15
$ws.Range("H5").Value = 17
$ws.Range("I5").Value = 13
$ws.Range("J5").Value = 13
$ws.Range("C6").Value = 24
$ws.Range("F6").Value = 23
$ws.Range("H6").Value = 21
$ws.Range("I6").Value = 18
$ws.Range("J6").Value = 26

# Garfield Park
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("F4").Value = 1
$ws.Range("E5").Value = 14
$ws.Range("E6").Value = 21
$ws.Range("F6").Value = 29

# Grand Crossing
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("C5").Value = 13
$ws.Range("E5").Value = 10
$ws.Range("I5").Value = 12
$ws.Range("C6").Value = 15
$ws.Range("E6").Value = 13
$ws.Range("I6").Value = 22

# Armour Square
$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("F5").Value = 4
$ws.Range("F6").Value = 5

# Uptown
$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("C5").Value = 4
$ws.Range("C6").Value = 5

# Englewood
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("B5").Value = 20
$ws.Range("G5").Value = 17
$ws.Range("B6").Value = 20
$ws.Range("G6").Value = 24

# Loop
$ws = $wb.Worksheets.Item("Loop")
$ws.Range("D2").Value = 6
$ws.Range("D6").Value = 23
$ws.Range("E6").Value = 36
$ws.Range("H6").Value = 26
$ws.Range("D7").Value = 42
$ws.Range("E7").Value = 47
$ws.Range("H7").Value = 36

# West Loop
$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("F6").Value = 4
$ws.Range("F7").Value = 5

# Printers Row
$ws = $wb.Worksheets.Item("Printers Row")
$ws.Range("D4").Value = 2
$ws.Range("D5").Value = 2

# Calumet Heights
$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("D3").Value = 2
$ws.Range("D5").Value = 2

# Edgewater
$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("D4").Value = 3
$ws.Range("D5").Value = 4

# Bridgeport
$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("B3").Value = 3
$ws.Range("B4").Value = 4

# Jefferson Park
$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("H4").Value = 3
$ws.Range("H5").Value = 3
